{"js": "// Applies the \"Monopoly Megaways\" copy refresh described in the commit\n// \"Added many more features\".\n//\n// Strategy: use Body.search() to locate the exact text runs that need\n// new wording and replace them in place with Range.insertText(..., \"Replace\").\n// This only rewrites the <w:t> text node of the matched run and leaves\n// sibling runs (including the leading empty <w:r/> anchor run used\n// throughout this document's bullet paragraphs) untouched.\n\nconst body = context.document.body;\n\nasync function replaceOnce(oldText, newText) {\n  const results = body.search(oldText, { matchCase: true });\n  results.load(\"items\");\n  await context.sync();\n  if (results.items.length === 0) {\n    throw new Error(\"Text not found: \" + oldText);\n  }\n  results.items[0].insertText(newText, \"Replace\");\n  await context.sync();\n}\n\n// 1) Main H1 title.\nawait replaceOnce(\n  \"Play Monopoly Megaways Free | 117,649 Ways to Win\",\n  \"Play Monopoly Megaways for Free\"\n);\n\n// 2) \"What we like\" bullets.\n// Remove the first bullet entirely (\"117,649 ways to win\").\nconst paras = body.paragraphs;\nparas.load(\"text\");\nawait context.sync();\n\nlet targetIndex = -1;\nfor (let i = 0; i < paras.items.length; i++) {\n  if (paras.items[i].text === \"117,649 ways to win\") {\n    targetIndex = i;\n    break;\n  }\n}\nif (targetIndex === -1) {\n  throw new Error(\"Could not find the '117,649 ways to win' bullet\");\n}\nparas.items[targetIndex].delete();\nawait context.sync();\n\n// Reword the \"highly volatile\" bullet.\nawait replaceOnce(\n  \"Highly volatile, with up to 14,700x your bet to win on free spins\",\n  \"Highly volatile, with big win potential\"\n);\n\n// Reword the \"Playable on all devices\" bullet, then add a brand-new bullet\n// right after it (\"Playable on all devices, including mobile\"), preserving\n// the paragraph's ListBullet formatting and its leading empty run.\nawait replaceOnce(\n  \"Playable on all devices\",\n  \"Mr. Monopoly adds extra features and increased multipliers\"\n);\n\nconst likeParas = body.paragraphs;\nlikeParas.load(\"text\");\nawait context.sync();\n\nlet mrMonopolyIndex = -1;\nfor (let i = 0; i < likeParas.items.length; i++) {\n  if (\n    likeParas.items[i].text ===\n    \"Mr. Monopoly adds extra features and increased multipliers\"\n  ) {\n    mrMonopolyIndex = i;\n    break;\n  }\n}\nif (mrMonopolyIndex === -1) {\n  throw new Error(\"Could not find the 'Mr. Monopoly' bullet\");\n}\n\nconst anchorPara = likeParas.items[mrMonopolyIndex];\nconst anchorRange = anchorPara.getRange();\nconst newBulletOoxml =\n  '<?xml version=\"1.0\" encoding=\"UTF-8\" standalone=\"yes\"?>' +\n  '<pkg:package xmlns:pkg=\"http://schemas.microsoft.com/office/2006/xmlPackage\">' +\n  '<pkg:part pkg:name=\"/word/document.xml\" pkg:contentType=\"application/vnd.openxmlformats-officedocument.wordprocessingml.document.main+xml\">' +\n  \"<pkg:xmlData>\" +\n  '<w:document xmlns:w=\"http://schemas.openxmlformats.org/wordprocessingml/2006/main\"><w:body>' +\n  '<w:p><w:pPr><w:pStyle w:val=\"ListBullet\"/><w:spacing w:line=\"240\" w:lineRule=\"auto\"/><w:ind w:left=\"720\"/></w:pPr><w:r/><w:r><w:t>Playable on all devices, including mobile</w:t></w:r></w:p>' +\n  \"</w:body></w:document>\" +\n  \"</pkg:xmlData></pkg:part></pkg:package>\";\nanchorRange.insertOoxml(newBulletOoxml, \"After\");\nawait context.sync();\n\n// 3) \"What we don't like\" bullets.\nawait replaceOnce(\n  \"May not be suitable for players who prefer low volatility slots\",\n  \"Limited variety of symbols\"\n);\nawait replaceOnce(\n  \"No progressive jackpot available\",\n  \"High-stakes gameplay may not be suitable for all players\"\n);\n\n// 4) Bold CTA paragraph (duplicate of the H1 title text).\nawait replaceOnce(\n  \"Play Monopoly Megaways Free | 117,649 Ways to Win\",\n  \"Play Monopoly Megaways for Free\"\n);\n\n// 5) Italic meta-description paragraph.\nawait replaceOnce(\n  \"Read our review of Monopoly Megaways, a highly volatile slot game with 117,649 ways to win from Big Time Gaming. Play for free on desktop or mobile.\",\n  \"Read our review of Monopoly Megaways and play this highly volatile slot game for free.\"\n);\n", "ps1": "# Applies the \"Monopoly Megaways\" copy refresh described in the commit\n# \"Added many more features\".\n#\n# Strategy: locate each paragraph that needs new wording by its current\n# (pre-edit) exact text, then rewrite that paragraph's Range via\n# Range.InsertXML() using a crafted <w:p> fragment. This only swaps the\n# paragraph's content while preserving its pPr (style/formatting) and the\n# leading empty <w:r/> anchor run pattern used throughout this document's\n# body/bullet paragraphs. A brand-new bullet paragraph is added the same\n# way, via InsertParagraphBefore() + InsertXML() on the freshly created\n# empty paragraph.\n\n$d = $word.ActiveDocument\n\nfunction Get-ParaIndexByText($doc, $text) {\n  $paras = $doc.Paragraphs\n  $count = $paras.Count\n  for ($i = 1; $i -le $count; $i++) {\n    $p = $paras.Item($i)\n    $t = $p.Range.Text\n    $t = $t.TrimEnd([char]13, [char]7)\n    if ($t -eq $text) {\n      return $i\n    }\n  }\n  return -1\n}\n\nfunction Set-ParagraphXml($doc, $index, $innerParaXml) {\n  $p = $doc.Paragraphs.Item($index)\n  $ooxml = '<?xml version=\"1.0\" encoding=\"UTF-8\" standalone=\"yes\"?>' + `\n    '<pkg:package xmlns:pkg=\"http://schemas.microsoft.com/office/2006/xmlPackage\">' + `\n    '<pkg:part pkg:name=\"/word/document.xml\" pkg:contentType=\"application/vnd.openxmlformats-officedocument.wordprocessingml.document.main+xml\">' + `\n    '<pkg:xmlData>' + `\n    '<w:document xmlns:w=\"http://schemas.openxmlformats.org/wordprocessingml/2006/main\"><w:body>' + `\n    $innerParaXml + `\n    '</w:body></w:document>' + `\n    '</pkg:xmlData></pkg:part></pkg:package>'\n  $p.Range.InsertXML($ooxml)\n}\n\nfunction Escape-Xml($text) {\n  return $text.Replace(\"&\", \"&amp;\").Replace(\"<\", \"&lt;\").Replace(\">\", \"&gt;\")\n}\n\nfunction Set-ParagraphTextByMatch($doc, $oldText, $newText, $innerXmlTemplate) {\n  $idx = Get-ParaIndexByText $doc $oldText\n  if ($idx -eq -1) {\n    throw \"Could not find paragraph with text: $oldText\"\n  }\n  $xml = $innerXmlTemplate -replace '__TEXT__', (Escape-Xml $newText)\n  Set-ParagraphXml $doc $idx $xml\n}\n\n$bulletTemplate = '<w:p><w:pPr><w:pStyle w:val=\"ListBullet\"/><w:spacing w:line=\"240\" w:lineRule=\"auto\"/><w:ind w:left=\"720\"/></w:pPr><w:r/><w:r><w:t>__TEXT__</w:t></w:r></w:p>'\n$headingTemplate = '<w:p><w:pPr><w:pStyle w:val=\"Heading1\"/></w:pPr><w:r><w:t>__TEXT__</w:t></w:r></w:p>'\n$boldTemplate = '<w:p><w:r/><w:r><w:rPr><w:b/></w:rPr><w:t>__TEXT__</w:t></w:r></w:p>'\n$italicTemplate = '<w:p><w:r/><w:r><w:rPr><w:i/></w:rPr><w:t>__TEXT__</w:t></w:r></w:p>'\n\n# 1) Main H1 title.\nSet-ParagraphTextByMatch $d \"Play Monopoly Megaways Free | 117,649 Ways to Win\" \"Play Monopoly Megaways for Free\" $headingTemplate\n\n# 2) \"What we like\" bullets.\n# Remove the first bullet entirely (\"117,649 ways to win\").\n$removeIdx = Get-ParaIndexByText $d \"117,649 ways to win\"\nif ($removeIdx -eq -1) {\n  throw \"Could not find the '117,649 ways to win' bullet\"\n}\n$d.Paragraphs.Item($removeIdx).Range.Delete()\n\n# Reword the \"highly volatile\" bullet.\nSet-ParagraphTextByMatch $d \"Highly volatile, with up to 14,700x your bet to win on free spins\" \"Highly volatile, with big win potential\" $bulletTemplate\n\n# Reword the \"Playable on all devices\" bullet.\nSet-ParagraphTextByMatch $d \"Playable on all devices\" \"Mr. Monopoly adds extra features and increased multipliers\" $bulletTemplate\n\n# Insert a brand-new bullet right after it.\n$anchorIdx = Get-ParaIndexByText $d \"Mr. Monopoly adds extra features and increased multipliers\"\nif ($anchorIdx -eq -1) {\n  throw \"Could not find the 'Mr. Monopoly' bullet\"\n}\n$nextPara = $d.Paragraphs.Item($anchorIdx + 1)\n$nextPara.Range.InsertParagraphBefore()\n$newBulletIdx = $anchorIdx + 1\n$newBulletXml = $bulletTemplate -replace '__TEXT__', (Escape-Xml \"Playable on all devices, including mobile\")\nSet-ParagraphXml $d $newBulletIdx $newBulletXml\n\n# 3) \"What we don't like\" bullets.\nSet-ParagraphTextByMatch $d \"May not be suitable for players who prefer low volatility slots\" \"Limited variety of symbols\" $bulletTemplate\nSet-ParagraphTextByMatch $d \"No progressive jackpot available\" \"High-stakes gameplay may not be suitable for all players\" $bulletTemplate\n\n# 4) Bold CTA paragraph (duplicate of the H1 title text).\nSet-ParagraphTextByMatch $d \"Play Monopoly Megaways Free | 117,649 Ways to Win\" \"Play Monopoly Megaways for Free\" $boldTemplate\n\n# 5) Italic meta-description paragraph.\nSet-ParagraphTextByMatch $d \"Read our review of Monopoly Megaways, a highly volatile slot game with 117,649 ways to win from Big Time Gaming. Play for free on desktop or mobile.\" \"Read our review of Monopoly Megaways and play this highly volatile slot game for free.\" $italicTemplate\n"}
